$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (prices and 1h volume %) per diff.
# Price values that parse as plain numbers need to be forced to stay
# text (matching the source inline-string cells with no special number
# formatting) by temporarily applying a Text format, then resetting the
# cell style back to Normal so no stray style/quote-prefix is left behind.

$ws.Range('D2').Value = '30.751.37'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').Value = '2.115.75'
$ws.Range('E3').Value = '  +10.33%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '333.73'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Range('E5').Value = '  +4.29%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5235'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Range('E7').Value = '  +3.54%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.4411'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Range('E8').Value = '  +8.34%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.09068'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Range('E9').Value = '  +8.69%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '46.64'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Range('E10').Value = '  +10.21%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.185'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Range('E11').Value = '  +6.80%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '25.27'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Range('E12').Value = '  +5.27%  '
$ws.Range('D13').Value = '2.117.43'
$ws.Range('E13').Value = '  +10.59%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.788'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Range('E14').Value = '  +5.72%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.784'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Range('E15').Value = '  +7.50%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '98.31'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Range('E16').Value = '  +6.26%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '1.001'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.00001138'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Range('E18').Value = '  +3.95%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.06643'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Range('E19').Value = '  +2.12%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '19.22'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Range('E20').Value = '  +3.81%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.404'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Range('E22').Value = '  +7.73%  '
$ws.Range('D23').Value = '30.856.25'
$ws.Range('E23').Value = '  +2.41%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '12.06'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Range('E24').Value = '  +6.13%  '
$ws.Range('D25').Value = '2.365.64'
$ws.Range('E25').Value = '  +10.86%  '
$ws.Range('E26').Value = '  +2.86%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '22.97'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Range('E27').Value = '  +4.75%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.556'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Range('E28').Value = '  +12.61%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '163.54'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '133.70'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Range('E30').Value = '  +3.70%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.190'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Range('E31').Value = '  +4.86%  '
$ws.Range('E32').Value = '  +2.37%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '6.254'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Range('E33').Value = '  +5.12%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '3.921'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.531'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Range('E35').Value = '  +27.85%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.02599'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Range('E36').Value = '  +5.84%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.595'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Range('E37').Value = '  +4.80%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '9.631'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Range('E38').Value = '  +12.00%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.06779'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Range('E39').Value = '  +5.21%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '12.82'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Range('E40').Value = '  +12.33%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.2276'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.6824'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Range('E42').Value = '  +4.84%  '
$ws.Range('E43').Value = '  +3.90%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '14.14'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Range('E44').Value = '  +5.66%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.9999'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.6392'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Range('E46').Value = '  +5.45%  '
$ws.Range('E47').Value = '  +3.19%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.674'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Range('E48').Value = '  +1.46%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.287'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Range('E49').Value = '  +6.26%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '83.22'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Range('E50').Value = '  +5.33%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.07080'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Range('E51').Value = '  +3.80%  '
